$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 9573.857
$ws.Range("J111").Value = 10000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -36134
$ws.Range("H137").Value = 3593.074
$ws.Range("I137").Value = 3418.0588
$ws.Range("J137").Value = 3890.6
$ws.Range("K137").Value = 10254.1764
$ws.Range("L137").Value = 11671.8
$ws.Range("M137").Value = -7704.1764
$ws.Range("N137").Value = -16771.8
$ws.Range("H138").Value = 2196.2292
$ws.Range("I138").Value = 1709.5172
$ws.Range("J138").Value = 2939.1052
$ws.Range("K138").Value = 5128.5516
$ws.Range("L138").Value = 8817.3156
$ws.Range("M138").Value = 11.44840000000022
$ws.Range("N138").Value = -19097.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1335.5
$ws.Range("I2").Value = 1566.6666
$ws.Range("J2").Value = 1104.3334
$ws.Range("K2").Value = 1566.6666
$ws.Range("L2").Value = 1104.3334
$ws.Range("M2").Value = -1453.6666
$ws.Range("N2").Value = -1330.3334
$ws.Range("H32").Value = 695314.1
$ws.Range("I32").Value = 781978.75
$ws.Range("J32").Value = 30885.666
$ws.Range("K32").Value = 781978.75
$ws.Range("L32").Value = 30885.666
$ws.Range("M32").Value = -781691.75
$ws.Range("N32").Value = -31459.666
$ws.Range("H45").Value = 3414.1428
$ws.Range("I45").Value = 3579.8
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 3579.8
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -3202.8
$ws.Range("N45").Value = -3754
$ws.Range("H61").Value = 2040.2325
$ws.Range("I61").Value = 1692.5
$ws.Range("K61").Value = 1692.5
$ws.Range("M61").Value = -1480.5
$ws.Range("H74").Value = 1127.8286
$ws.Range("I74").Value = 786.8889
$ws.Range("J74").Value = 2278.5
$ws.Range("K74").Value = 786.8889
$ws.Range("L74").Value = 2278.5
$ws.Range("M74").Value = 87.11109999999996
$ws.Range("N74").Value = -4026.5
$ws.Range("H77").Value = 1127.8286
$ws.Range("I77").Value = 786.8889
$ws.Range("J77").Value = 2278.5
$ws.Range("K77").Value = 3934.4445
$ws.Range("L77").Value = 11392.5
$ws.Range("M77").Value = 433.5554999999999
$ws.Range("N77").Value = -20128.5
$ws.Range("H97").Value = 971.7143
$ws.Range("J97").Value = 1156
$ws.Range("L97").Value = 1156
$ws.Range("N97").Value = -2148
$ws.Range("H116").Value = 1335.5
$ws.Range("I116").Value = 1566.6666
$ws.Range("J116").Value = 1104.3334
$ws.Range("K116").Value = 1566.6666
$ws.Range("L116").Value = 1104.3334
$ws.Range("M116").Value = 727.3334
$ws.Range("N116").Value = -5692.3334
$ws.Range("H132").Value = 4294.3145
$ws.Range("I132").Value = 2818.318
$ws.Range("K132").Value = 8454.954000000002
$ws.Range("M132").Value = -5924.954000000002
$ws.Range("H136").Value = 2040.2325
$ws.Range("I136").Value = 1692.5
$ws.Range("K136").Value = 5077.5
$ws.Range("M136").Value = -2527.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1335.5
$ws.Range("I3").Value = 1566.6666
$ws.Range("J3").Value = 1104.3334
$ws.Range("K3").Value = 1566.6666
$ws.Range("L3").Value = 1104.3334
$ws.Range("M3").Value = -1452.6666
$ws.Range("N3").Value = -1332.3334
$ws.Range("H94").Value = 1641.8
$ws.Range("I94").Value = 1369.6666
$ws.Range("J94").Value = 2050
$ws.Range("K94").Value = 1369.6666
$ws.Range("L94").Value = 2050
$ws.Range("M94").Value = -918.6666
$ws.Range("N94").Value = -2952
$ws.Range("H107").Value = 63756.875
$ws.Range("I107").Value = 67807.336
$ws.Range("K107").Value = 67807.336
$ws.Range("M107").Value = -65887.336
$ws.Range("H132").Value = 40524.117
$ws.Range("J132").Value = 40524.117
$ws.Range("L132").Value = 40524.117
$ws.Range("N132").Value = -50644.117
$ws.Range("H134").Value = 2325.0833
$ws.Range("I134").Value = 1705.6923
$ws.Range("K134").Value = 5117.0769
$ws.Range("M134").Value = -2582.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 60333.332
$ws.Range("J28").Value = 60333.332
$ws.Range("L28").Value = 60333.332
$ws.Range("N28").Value = -60823.332
$ws.Range("H31").Value = 6871.1523
$ws.Range("I31").Value = 1429.4762
$ws.Range("J31").Value = 11442.16
$ws.Range("K31").Value = 1429.4762
$ws.Range("L31").Value = 11442.16
$ws.Range("M31").Value = -1134.4762
$ws.Range("N31").Value = -12032.16
$ws.Range("H34").Value = 6871.1523
$ws.Range("I34").Value = 1429.4762
$ws.Range("J34").Value = 11442.16
$ws.Range("K34").Value = 1429.4762
$ws.Range("L34").Value = 11442.16
$ws.Range("M34").Value = -1227.4762
$ws.Range("N34").Value = -11846.16
$ws.Range("H93").Value = 25161.4
$ws.Range("I93").Value = 5203.5
$ws.Range("J93").Value = 38466.668
$ws.Range("K93").Value = 5203.5
$ws.Range("L93").Value = 38466.668
$ws.Range("M93").Value = -3331.5
$ws.Range("N93").Value = -42210.668
$ws.Range("H95").Value = 18999.666
$ws.Range("J95").Value = 18999.666
$ws.Range("L95").Value = 18999.666
$ws.Range("N95").Value = -24491.666
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H105").Value = 2124.875
$ws.Range("I105").Value = 1800
$ws.Range("K105").Value = 1800
$ws.Range("M105").Value = -53
$ws.Range("H107").Value = 824.1
$ws.Range("I107").Value = 328.2
$ws.Range("J107").Value = 1320
$ws.Range("K107").Value = 328.2
$ws.Range("L107").Value = 1320
$ws.Range("M107").Value = 1591.8
$ws.Range("N107").Value = -5160
$ws.Range("H122").Value = 1925.641
$ws.Range("I122").Value = 1883.8889
$ws.Range("J122").Value = 1961.4286
$ws.Range("K122").Value = 5651.6667
$ws.Range("L122").Value = 5884.2858
$ws.Range("M122").Value = -3201.6667
$ws.Range("N122").Value = -10784.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 739.0476
$ws.Range("J5").Value = 1276.8334
$ws.Range("L5").Value = 3830.5002
$ws.Range("N5").Value = -4054.5002
$ws.Range("H109").Value = 3478.2354
$ws.Range("I109").Value = 822.2222
$ws.Range("J109").Value = 6466.25
$ws.Range("K109").Value = 2466.6666
$ws.Range("L109").Value = 19398.75
$ws.Range("M109").Value = -1426.6666
$ws.Range("N109").Value = -21478.75
$ws.Range("H129").Value = 2002.2142
$ws.Range("I129").Value = 470
$ws.Range("J129").Value = 2615.1
$ws.Range("K129").Value = 1410
$ws.Range("L129").Value = 7845.299999999999
$ws.Range("M129").Value = 3590
$ws.Range("N129").Value = -17845.3
$ws.Range("H131").Value = 939.5217
$ws.Range("J131").Value = 1207.6923
$ws.Range("L131").Value = 3623.0769
$ws.Range("N131").Value = -13703.0769
$ws.Range("H135").Value = 739.0476
$ws.Range("J135").Value = 1276.8334
$ws.Range("L135").Value = 11491.5006
$ws.Range("N135").Value = -16561.5006
$ws.Range("H139").Value = 5049.6875
$ws.Range("I139").Value = 2163
$ws.Range("K139").Value = 6489
$ws.Range("M139").Value = -1349
$ws.Range("H140").Value = 1839.5454
$ws.Range("I140").Value = 1756.25
$ws.Range("K140").Value = 5268.75
$ws.Range("M140").Value = -88.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 39990
$ws.Range("J51").Value = 39990
$ws.Range("L51").Value = 39990
$ws.Range("N51").Value = -41008
$ws.Range("H99").Value = 20117.75
$ws.Range("I99").Value = 15235.5
$ws.Range("J99").Value = 25000
$ws.Range("K99").Value = 15235.5
$ws.Range("L99").Value = 25000
$ws.Range("M99").Value = -12989.5
$ws.Range("N99").Value = -29492
$ws.Range("H113").Value = 251475.25
$ws.Range("I113").Value = 500444
$ws.Range("J113").Value = 2506.5
$ws.Range("K113").Value = 500444
$ws.Range("L113").Value = 2506.5
$ws.Range("M113").Value = -498274
$ws.Range("N113").Value = -6846.5
$ws.Range("H132").Value = 3514.0833
$ws.Range("I132").Value = 2796
$ws.Range("J132").Value = 4362.727
$ws.Range("K132").Value = 8388
$ws.Range("L132").Value = 13088.181
$ws.Range("M132").Value = -5858
$ws.Range("N132").Value = -18148.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2060.7334
$ws.Range("I132").Value = 1338.5217
$ws.Range("K132").Value = 4015.5651
$ws.Range("M132").Value = -1485.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5382.857
$ws.Range("I96").Value = 4020
$ws.Range("J96").Value = 7200
$ws.Range("K96").Value = 4020
$ws.Range("L96").Value = 7200
$ws.Range("M96").Value = -2647
$ws.Range("N96").Value = -9946
$ws.Range("H100").Value = 1295.8948
$ws.Range("I100").Value = 1524.3334
$ws.Range("J100").Value = 904.2857
$ws.Range("K100").Value = 3048.6668
$ws.Range("L100").Value = 1808.5714
$ws.Range("M100").Value = -2507.6668
$ws.Range("N100").Value = -2890.5714
$ws.Range("H113").Value = 1054.909
$ws.Range("I113").Value = 1256.3529
$ws.Range("J113").Value = 370
$ws.Range("K113").Value = 3769.0587
$ws.Range("L113").Value = 1110
$ws.Range("M113").Value = -1599.0587
$ws.Range("N113").Value = -5450
